# Update variable_type column (D) to "list" for the rows identified in the
# commit "Work through all vars". Rows 65 and 67 are already "list" and are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(53, 54, 55, 56, 57, 58, 59, 60, 61, 62, 63, 64, 66, 68, 69, 70, 71, 72)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "list"
}
